# ===========================================================================
# Allure fix and Digital assessment new query implemented
#
# - Adds a new "GradeOneStudentCredentials" worksheet (3rd tab) with a new
#   set of generated grade-1 student credentials.
# - Tweaks column widths / selections on the two existing sheets.
# ===========================================================================

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "GradeOneStudentCredentials" as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "GradeOneStudentCredentials"

# Header row (same headers as the other credential sheets)
$ws3.Range("A1").Value = "UserId"
$ws3.Range("B1").Value = "Password"
$ws3.Range("C1").Value = "UserName"
$ws3.Range("D1").Value = "Signature"

# Header cell styles (match the "Normal 2" cell style used on the other sheets)
$ws3.Range("A1").Style = "Normal 2"
$ws3.Range("B1").Style = "Normal 2"

# Data rows - newly generated grade-one student logins
$ws3.Range("B2").Value = "Password@123"
$ws3.Range("A2").Value = "AutoGoHlE_208"

$ws3.Range("A3").Value = "AutoVPGad_263"

$ws3.Range("A4").Value = "AutoglwqZ_839"
$ws3.Range("B4").Value = "Password@123"

$ws3.Range("A5").Value = "AutoiEIjf_985"
$ws3.Range("B5").Value = "Password@123"

# Column widths for the new sheet
$ws3.Columns.Item(1).ColumnWidth = 14.5
$ws3.Columns.Item(2).ColumnWidth = 13.5

# Page setup (matches the other credential sheets)
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Selection / active cell on the new sheet
$ws3.Range("A6").Select()

# --- ParentCredentials (sheet1) view + column tweaks ---
$ws1 = $wb.Worksheets.Item("ParentCredentials")
$ws1.Columns.Item(4).ColumnWidth = 10.333333333333332
$ws1.Range("G10").Select()

# --- StudentCredentials (sheet2) view + column tweaks ---
$ws2 = $wb.Worksheets.Item("StudentCredentials")
$ws2.Columns.Item(1).ColumnWidth = 7
$ws2.Columns.Item(2).ColumnWidth = 8.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 8
$ws2.Columns.Item(4).ColumnWidth = 7.500000000000001
$ws2.Range("A1:D2").Select()

# Make the new sheet the active / visible tab (tabSelected moves here)
$ws3.Activate()

# Maximize the workbook window
$excel.ActiveWindow.WindowState = -4143
